# Burn Down Chart update
# - Extends the burn-down data from 18 rows (Sheet1!A2:C19) to 26 rows (Sheet1!A2:C27)
#   so that every calendar day (including weekends) between Feb 27 and Mar 24 2015
#   is represented, instead of only weekdays.
# - Updates the "Hours Left" (col B) and "Burn-Down" (col C) figures to match.
# - The line chart series (which reference Sheet1!A2:A19 / B / C) automatically
#   grow to the new range once the worksheet data & formulas are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure the number formats used by the existing table extend across the
# full new range (date format for column A, one-decimal numeric format for B/C)
# so that newly-created rows 20:27 pick up the same look as the existing rows.
$ws.Range("A2:A27").NumberFormat = "mmmm\ d"
$ws.Range("B2:C27").NumberFormat = "0.0"

# --- Row 2 (2015-02-27) ---
$ws.Range("A2").Value = 42062
$ws.Range("B2").Formula = "=B3+10.8"
$ws.Range("C2").Value = 128

# --- Row 3 (2015-02-28) - new row ---
$ws.Range("A3").Value = 42063
$ws.Range("B3").Formula = "=B4+10.8"
$ws.Range("C3").Formula = "=C2"

# --- Row 4 (2015-03-01) - new row ---
$ws.Range("A4").Value = 42064
$ws.Range("B4").Formula = "=B5+10.8"
$ws.Range("C4").Formula = "=C3"

# --- Row 5 (2015-03-02) ---
$ws.Range("A5").Value = 42065
$ws.Range("B5").Formula = "=B6+10.8"
$ws.Range("C5").Formula = "=C2-2"

# --- Row 6 (2015-03-03) ---
$ws.Range("A6").Value = 42066
$ws.Range("B6").Formula = "=B7+10.8"
$ws.Range("C6").Formula = "=C5-2"

# --- Row 7 (2015-03-04) ---
$ws.Range("A7").Value = 42067
$ws.Range("B7").Formula = "=B8+10.8"
$ws.Range("C7").Formula = "=C6-124+112.5"

# --- Row 8 (2015-03-05) ---
$ws.Range("A8").Value = 42068
$ws.Range("B8").Formula = "=B9+10.8"
$ws.Range("C8").Formula = "=C7"

# --- Row 9 (2015-03-06) ---
$ws.Range("A9").Value = 42069
$ws.Range("B9").Formula = "=B10+10.8"
$ws.Range("C9").Formula = "=C8"

# --- Row 10 (2015-03-07) - new row ---
$ws.Range("A10").Value = 42070
$ws.Range("B10").Formula = "=B11+10.8"
$ws.Range("C10").Formula = "=C9"

# --- Row 11 (2015-03-08) - new row ---
$ws.Range("A11").Value = 42071
$ws.Range("B11").Formula = "=B12+10.8"
$ws.Range("C11").Formula = "=C10"

# --- Row 12 (2015-03-09) ---
$ws.Range("A12").Value = 42072
$ws.Range("B12").Formula = "=B13+10.8"
$ws.Range("C12").Formula = "=C9"

# --- Row 13 (2015-03-10) ---
$ws.Range("A13").Value = 42073
$ws.Range("B13").Formula = "=B14+10.8"
$ws.Range("C13").Formula = "=C12"

# --- Row 14 (2015-03-11) ---
$ws.Range("A14").Value = 42074
$ws.Range("B14").Formula = "=B15+10.8"
$ws.Range("C14").Formula = "=C13"

# --- Row 15 (2015-03-12) ---
$ws.Range("A15").Value = 42075
$ws.Range("B15").Formula = "=B16+10.8"
$ws.Range("C15").Formula = "=C14"

# --- Row 16 (2015-03-13) ---
$ws.Range("A16").Value = 42076
$ws.Range("B16").Formula = "=B17+10.8"
$ws.Range("C16").Formula = "=C15-5"

# --- Row 17 (2015-03-14) - new row ---
$ws.Range("A17").Value = 42077
$ws.Range("B17").Formula = "=B18+10.8"
$ws.Range("C17").Formula = "=C16"

# --- Row 18 (2015-03-15) - new row ---
$ws.Range("A18").Value = 42078
$ws.Range("B18").Formula = "=B19+10.8"
$ws.Range("C18").Formula = "=C17"

# --- Row 19 (2015-03-16) ---
$ws.Range("A19").Value = 42079
$ws.Range("B19").Formula = "=B20+10.8"
$ws.Range("C19").Formula = "=C16"

# --- Row 20 (2015-03-17) ---
$ws.Range("A20").Value = 42080
$ws.Range("B20").Formula = "=B21+10.8"
$ws.Range("C20").Formula = "=C19"

# --- Row 21 (2015-03-18) ---
$ws.Range("A21").Value = 42081
$ws.Range("B21").Formula = "=B22+10.8"
$ws.Range("C21").Formula = "=C20"

# --- Row 22 (2015-03-19) ---
$ws.Range("A22").Value = 42082
$ws.Range("B22").Formula = "=B23+10.8"
$ws.Range("C22").Formula = "=C21"

# --- Row 23 (2015-03-20) ---
$ws.Range("A23").Value = 42083
$ws.Range("B23").Formula = "=B24+10.8"
$ws.Range("C23").Formula = "=C22"

# --- Row 24 (2015-03-21) - new row ---
$ws.Range("A24").Value = 42084
$ws.Range("B24").Formula = "=B25+10.8"
$ws.Range("C24").Formula = "=C23"

# --- Row 25 (2015-03-22) - new row ---
$ws.Range("A25").Value = 42085
$ws.Range("B25").Formula = "=B26+10.8"
$ws.Range("C25").Formula = "=C24"

# --- Row 26 (2015-03-23) ---
$ws.Range("A26").Value = 42086
$ws.Range("B26").Formula = "=B27+10.8"
$ws.Range("C26").Formula = "=C23"

# --- Row 27 (2015-03-24) ---
$ws.Range("A27").Value = 42087
$ws.Range("B27").Value = 0
$ws.Range("C27").Formula = "=C26"

# Match the selection left behind by the author when they finished editing.
[void]$ws.Range("D18").Select()
